$d = $word.ActiveDocument

$replacements = @(
    @{old = "561÷6="; new = "319÷9="},
    @{old = "174÷7="; new = "487÷9="},
    @{old = "544÷9="; new = "516÷5="},
    @{old = "493÷3="; new = "698÷2="},
    @{old = "562÷6="; new = "200÷7="},
    @{old = "195÷4="; new = "836÷9="},
    @{old = "729÷8="; new = "360÷9="},
    @{old = "783÷4="; new = "694÷8="},
    @{old = "777÷6="; new = "711÷8="},
    @{old = "101÷2="; new = "417÷5="},
    @{old = "942÷9="; new = "390÷9="},
    @{old = "426÷5="; new = "382÷2="},
    @{old = "349÷5="; new = "331÷2="},
    @{old = "861÷4="; new = "571÷5="},
    @{old = "584÷4="; new = "595÷9="},
    @{old = "309÷5="; new = "131÷2="},
    @{old = "621÷2="; new = "433÷7="},
    @{old = "237÷8="; new = "676÷8="},
    @{old = "469÷4="; new = "566÷7="},
    @{old = "465÷2="; new = "725÷9="},
    @{old = "379÷3="; new = "203÷4="},
    @{old = "598÷5="; new = "287÷6="},
    @{old = "285÷7="; new = "616÷8="},
    @{old = "116÷2="; new = "482÷2="},
    @{old = "305÷2="; new = "657÷7="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
